# Add a second data-driven test sheet ("calculateCarLoan") after the
# existing "calculateTax" sheet, with its own header row + sample data,
# and make it the active sheet (mirrors the author's "2nd Data Driven
# Test Case added" commit).

$wb = $excel.ActiveWorkbook

# Insert the new worksheet after the last existing sheet (i.e. after
# "calculateTax") so it lands as the 2nd tab.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$carLoanSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$carLoanSheet.Name = "calculateCarLoan"

# Header row
$carLoanSheet.Range("A1").Value = "Loan_Amount"
$carLoanSheet.Range("B1").Value = "Loan_Period"
$carLoanSheet.Range("C1").Value = "EMI_Starts_From"
$carLoanSheet.Range("D1").Value = "Interest_Rate"
$carLoanSheet.Range("E1").Value = "Upfront_Charges"

# Data row 2
$carLoanSheet.Range("A2").Value = 2000000
$carLoanSheet.Range("B2").Value = 10
$carLoanSheet.Range("C2").Value = "At the time of loan disbursement"
$carLoanSheet.Range("D2").Value = 6.5
$carLoanSheet.Range("E2").Value = 100000

# Data row 3
$carLoanSheet.Range("A3").Value = 1500000
$carLoanSheet.Range("B3").Value = 15
$carLoanSheet.Range("C3").Value = "From next month after disbursement"
$carLoanSheet.Range("D3").Value = 7
$carLoanSheet.Range("E3").Value = 50000

# Make the newly added sheet the active/selected tab, like in the diff
# (tabSelected moves from "calculateTax" to "calculateCarLoan").
$carLoanSheet.Activate()
